$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''39.918.57'
$ws.Range("E2").Value = '  -4.18%  '
$ws.Range("D3").Value = '''2.344.68'
$ws.Range("E3").Value = '  -5.30%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''308.23'
$ws.Range("E5").Value = '  -3.75%  '
$ws.Range("D6").Value = '''84.05'
$ws.Range("E6").Value = '  -9.01%  '
$ws.Range("E7").Value = '  -3.81%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("E9").Value = '  -5.48%  '
$ws.Range("D10").Value = '''0.0819'
$ws.Range("E10").Value = '  -5.13%  '
$ws.Range("D11").Value = '''30.03'
$ws.Range("E11").Value = '  -9.05%  '
$ws.Range("E12").Value = '  -1.02%  '
$ws.Range("D13").Value = '''2.710.91'
$ws.Range("E13").Value = '  -5.04%  '
$ws.Range("D14").Value = '''6.40'
$ws.Range("E14").Value = '  -7.08%  '
$ws.Range("D15").Value = '''14.81'
$ws.Range("E15").Value = '  -4.53%  '
$ws.Range("D16").Value = '''2.360.64'
$ws.Range("E16").Value = '  -4.10%  '
$ws.Range("D17").Value = '''0.750'
$ws.Range("E17").Value = '  -5.71%  '
$ws.Range("D18").Value = '''39.978.40'
$ws.Range("E18").Value = '  -3.93%  '
$ws.Range("D19").Value = '''0.0₃0900'
$ws.Range("E19").Value = '  -4.46%  '
$ws.Range("D20").Value = '''6.07'
$ws.Range("E20").Value = '  -5.77%  '
$ws.Range("D21").Value = '''67.82'
$ws.Range("E21").Value = '  -3.92%  '
$ws.Range("D22").Value = '''10.63'
$ws.Range("E22").Value = '  -5.54%  '
$ws.Range("D23").Value = '''234.24'
$ws.Range("E23").Value = '  -2.38%  '
$ws.Range("D24").Value = '''2.55'
$ws.Range("E24").Value = '  -7.26%  '
$ws.Range("D25").Value = '''1.00'
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("E26").Value = '  -7.97%  '
$ws.Range("D27").Value = '''23.40'
$ws.Range("E27").Value = '  -6.50%  '
$ws.Range("E28").Value = '  -1.79%  '
$ws.Range("D29").Value = '''9.19'
$ws.Range("E29").Value = '  -5.79%  '
$ws.Range("D30").Value = '''34.23'
$ws.Range("E30").Value = '  -6.50%  '
$ws.Range("D31").Value = '''151.62'
$ws.Range("E31").Value = '  -3.49%  '
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("D33").Value = '''5.12'
$ws.Range("E33").Value = '  -5.77%  '
$ws.Range("D34").Value = '''0.0724'
$ws.Range("E34").Value = '  -5.23%  '
$ws.Range("E35").Value = '  -5.51%  '
$ws.Range("E36").Value = '  -2.42%  '
$ws.Range("D37").Value = '''2.76'
$ws.Range("E37").Value = '  -4.44%  '
$ws.Range("D38").Value = '''0.0989'
$ws.Range("E38").Value = '  -4.69%  '
$ws.Range("D39").Value = '''15.64'
$ws.Range("E39").Value = '  -8.94%  '
$ws.Range("D40").Value = '''1.69'
$ws.Range("E40").Value = '  -8.31%  '
$ws.Range("D41").Value = '''3.80'
$ws.Range("E41").Value = '  -5.35%  '
$ws.Range("D42").Value = '''2.36'
$ws.Range("E42").Value = '  -3.93%  '
$ws.Range("D43").Value = '''1.950.98'
$ws.Range("E43").Value = '  -2.59%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '''17.78'
$ws.Range("E44").Value = '  -4.94%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = '''0.0265'
$ws.Range("E45").Value = '  -6.87%  '
$ws.Range("D46").Value = '''9.35'
$ws.Range("E46").Value = '  -1.89%  '
$ws.Range("D47").Value = '''2.64'
$ws.Range("E47").Value = '  -11.11%  '
$ws.Range("D48").Value = '''2.587.00'
$ws.Range("E48").Value = '  -4.73%  '
$ws.Range("D49").Value = '''92.21'
$ws.Range("E49").Value = '  -5.56%  '
$ws.Range("D50").Value = '''70.80'
$ws.Range("E50").Value = '  -6.32%  '
$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D51").Value = '''63.47'
$ws.Range("E51").Value = '  -5.77%  '
